# Add a "preview" section (table + line preview queries) to the metadata
# sheet, inserted right after the dataset.status row (i.e. becoming the new
# rows 4 and 5, pushing the rest of the metadata rows down by two).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Build the multi-line "source(...)"-style preview formulas exactly as they
# appear in the shared strings table.
$tablePreview = @"
source(ds:'{{dataset.id}}');
query([
  { dim:'time', role:'row', items:[] },
  { dim:'indicator', role:'col', items:[] } 
]);
format(p:3);
order(dir:'row', index:-1, asc:'az');
limit(start:0, length:5);
"@

$linePreview = @"
source(ds:'{{dataset.id}}');
query([
  { dim:'time', role:'row', items:[] },
  { dim:'indicator', role:'col', items:[] } 
]);
format(p:3);
order(dir:'row', index:-1, asc:'az');
line(x:-1);
"@

# Insert two new blank rows at row 4 (the old "dataset.commit.id" row and
# everything below it shifts down by two rows).
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

# Fill in the new rows with the preview keys/values.
$ws.Range("A4").Value = "dataset.preview.table"
$ws.Range("B4").Value = $tablePreview
$ws.Range("A5").Value = "dataset.preview.line"
$ws.Range("B5").Value = $linePreview

# Match the plain "key/value" styling used elsewhere (vertical centering),
# plus turn on word wrap so the multi-line formulas are fully visible, and
# grow the rows to fit the wrapped text.
foreach ($addr in @("A4", "B4", "A5", "B5")) {
    $cell = $ws.Range($addr)
    $cell.VerticalAlignment = -4108
    $cell.WrapText = $true
    $cell.VerticalAlignment = -4108
}

$ws.Rows.Item(4).RowHeight = 120
$ws.Rows.Item(5).RowHeight = 120

# Restore the active selection to B7 (the "dataset.commit.HEAD" value cell
# after the insertion).
$ws.Range("B7").Select()
